$wb = $excel.ActiveWorkbook

# Sheet: 展览
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 456
$ws1.Range("F7").Value = 367
$ws1.Range("F8").Value = 2078
$ws1.Range("F16").Value = 1400
$ws1.Range("F19").Value = 525
$ws1.Range("F22").Value = 7141
$ws1.Range("F23").Value = 7774
$ws1.Range("F24").Value = 42
$ws1.Range("F28").Value = 88
$ws1.Range("F30").Value = 259
$ws1.Range("F44").Value = 330
$ws1.Range("F45").Value = 237
$ws1.Range("F47").Value = 83

# Sheet: 演出
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F11").Value = 16
$ws2.Range("F17").Value = 294

# Sheet: 本地生活
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F3").Value = 2603
$ws3.Range("F4").Value = 275

# Sheet: 全部类型
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 456
$ws4.Range("F4").Value = 275
$ws4.Range("F10").Value = 2078
$ws4.Range("F18").Value = 1400
$ws4.Range("F20").Value = 525
$ws4.Range("F22").Value = 7141
$ws4.Range("F23").Value = 7774
$ws4.Range("F24").Value = 42
$ws4.Range("F27").Value = 88
$ws4.Range("F28").Value = 259
$ws4.Range("F37").Value = 16
$ws4.Range("F42").Value = 330
$ws4.Range("F43").Value = 237
$ws4.Range("F45").Value = 83
$ws4.Range("F49").Value = 294

$wb.Save()
